# Duplicate the "T2" worksheet into a new "T3" worksheet (third tab),
# update the duplicated data row with the new subnet (192.168.75.0/29),
# attach fresh tables to the new sheet, and select/activate it - matching
# what Excel produces when you right-click a tab -> "Move or Copy" ->
# "Create a copy" and then edit the new sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("T2")

# Copy T2 and place the copy immediately after it -> becomes the 3rd tab.
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "T3"

# New subnet data (192.168.75.0/29) replacing the copied 192.168.47.0/28 row.
$ws3.Range("A2").Value = "192.168.75.0"
$ws3.Range("B2").Value = "255.255.255.248 ó /29"
$ws3.Range("C2").Value = "192.168.75.1"
$ws3.Range("D2").Value = "192.168.75.6"
$ws3.Range("E2").Value = "192.168.75.7"
$ws3.Range("F2").Value = 8
$ws3.Range("G2").Value = 7

# Re-create the two tables (lost on sheet copy) on the new sheet.
$lo1 = $ws3.ListObjects.Add(1, $ws3.Range("A1:G2"), $null, 1)
$lo1.Name = "Tabla252"
$lo1.TableStyle = "TableStyleDark7"

$lo2 = $ws3.ListObjects.Add(1, $ws3.Range("A4:E11"), $null, 1)
$lo2.Name = "Tabla367"
$lo2.TableStyle = "TableStyleDark6"

# Make T3 the active sheet/tab with the same selection the author left it in.
$ws3.Activate()
$ws3.Range("C16").Select() | Out-Null
